# Update post oliver Review
# - Rework the "ZeroRisika" Minus text (shorten it)
# - Replace "CaptainDelegate" Plus text with the delegation description
# - Replace the "SolveZero" hero with a new "FifthWhy" hero (name, plus, minus)
# - Replace "WikiLeakz" Plus text with documentation description
# - Replace "Hu-Dini" Minus text
# - Replace "Khon-Troller" Minus text
# - Replace "iRobot" Plus text
# - Reset the sheet scroll position and change the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B6").Value = "She never makes any mistake"

$ws.Range("B7").Value = "He delegates everything & always at the right person"

$ws.Range("A8").Value = "FifthWhy"
$ws.Range("B8").Value = 'She can make anyone search for the deepest "why"'
$ws.Range("D8").Value = "the team loses efficiency and commitment tends to go down the drain"

$ws.Range("B9").Value = "Everything he does is perfectly and extensively documented"

$ws.Range("D12").Value = "but the team avoids asking when they really need help"

$ws.Range("D14").Value = "emotions are ignored and metrics get gamed"

$ws.Range("B16").Value = "Removes emotion from every discussion"

# Reset scroll position back to A1 (was topLeftCell="C1") and move the
# selection to D2:D17 (was E26)
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("D2:D17").Select()
